$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MockClientMiddleware tests now working: the "GetB Method Authorization"
# scenarios previously expected a 404, but now correctly expect a 403
# (Forbidden) response.
$ws.Range("G9").Value = 403
$ws.Range("G15").Value = 403

# The refactored tests (distinct ports for each mock client API) made some
# of the sample JSON/test-case values longer, so those columns were widened
# by hand (losing their "best fit" auto-sizing) when the sheet was last
# saved.
$ws.Columns.Item(2).ColumnWidth = 15.451822916666666
$ws.Columns.Item(4).ColumnWidth = 34.307291666666664
$ws.Columns.Item(6).ColumnWidth = 22.736979166666668
$ws.Columns.Item(7).ColumnWidth = 8.736979166666666

# The active selection when the file was last saved was cell D8.
[void]$ws.Range("D8").Select()
